# Aplicados los cambios: "Agregados aerodromos hechos por Luis y actualizada hoja de calculo"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AD")

# Color verde (fillId con fgColor FF92D050) usado para marcar filas "Convertido = si"
$greenFill = 5296274   # BGR de RGB(92,D0,50)
$center = -4108         # xlCenter

# Filas convertidas por Luis (columna B ya contiene "Luis")
$luisRows = @(71,72,76,77,78,79,80,82,84,85,87,88,89,90,91)

# Filas convertidas por Jose (columna B ya contiene "Jose"), salvo la 200 -> Tomas
$joseRows = @(141,142,143,144,174,177,203,207,208)

foreach ($r in $luisRows) {
    $a = $ws.Cells.Item($r, 1)
    $a.Value2 = "si"
    $a.Interior.Color = $greenFill
    $a.HorizontalAlignment = $center

    $b = $ws.Cells.Item($r, 2)
    $b.Interior.Color = $greenFill
    $b.HorizontalAlignment = $center
}

foreach ($r in $joseRows) {
    $a = $ws.Cells.Item($r, 1)
    $a.Value2 = "si"
    $a.Interior.Color = $greenFill
    $a.HorizontalAlignment = $center

    $b = $ws.Cells.Item($r, 2)
    $b.Interior.Color = $greenFill
    $b.HorizontalAlignment = $center
}

# Fila 200: convertida, pero el aerodromo fue hecho por Tomas (no Jose)
$a200 = $ws.Cells.Item(200, 1)
$a200.Value2 = "si"
$a200.Interior.Color = $greenFill
$a200.HorizontalAlignment = $center

$b200 = $ws.Cells.Item(200, 2)
$b200.Value2 = "Tomas"
$b200.Interior.Color = $greenFill
$b200.HorizontalAlignment = $center

# Comentario en la fila 55: encontrado en xplane.es
$ws.Range("I55").Value2 = "En xplane.es"

# Actualizar la vista de la hoja (celda activa / seleccion)
$ws.Activate() | Out-Null
$ws.Range("I200").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 167
$win.ScrollColumn = 1
